$d = $word.ActiveDocument

# The document currently ends with a "Requisitos" Heading2 paragraph
# followed by a ListBullet paragraph listing the prerequisite course
# ("LOM3089 -  Mecânica dos Fluidos e Reologia  (Requisito fraco)").
# Both paragraphs must be removed so that the document ends right after
# the Bibliografia section.

$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Requisitos") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $range = $d.Range($target.Range.Start, $d.Content.End)
    $range.Delete()
}
